$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: % value changes from 20 to 70
$ws.Range("C15").Value = 70

# Row 16: new C16 value (0)
$ws.Range("C16").Value = 0

# Row 17: new row with C17 = 50 and D17 = "Validation"
$ws.Range("C17").Value = 50
$ws.Range("D17").Value = "Validation"

# Update the selection to reflect the new active cell (C18)
$ws.Range("C18").Select()
